# Fix error in live simulations from NAT vote share
# Update the raw poll figures in row 2 (Sheet1) that feed the normaliser.
# All downstream formulas (K2 SUMIF, row 4 percentages, row 6 IFNA values,
# and K8 SUMPRODUCT) recalculate automatically from these inputs.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

$ws.Range("A2").Value = 30
$ws.Range("B2").Value = 31
$ws.Range("C2").Value = 8
$ws.Range("D2").Value = 13
$ws.Range("I2").Value = 18

$excel.Calculate()

# Move active selection from L17 to K17 as in the saved workbook view.
$ws.Range("K17").Select()
